# Update workbook metadata and terminology references to match the
# published CDA FHIR logical model release (2.0.0-sd-202406-matchbox-patch).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet ---------------------------------------------------
$wsElements = $wb.Worksheets.Item("Elements")

# IVL_INT.operator binding value set
$wsElements.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"

# Column Z widened to fit the new (longer) value-set URL text.
$wsElements.Columns.Item(26).ColumnWidth = 50.3814697265625
